# Textbox response formatting fix
# Renames the task-order sheets and refreshes the timestamped stim/response
# filenames referenced in column B of each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511687042964354"
$ws1.Range("B2").Value = "go_stims-16511687042684026.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687042794397.csv"
$ws1.Range("B4").Value = "go_stims-16511687042804015.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687042954304.csv"

# --- Sheet 2: NB_TO ------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16511687083751464"
$ws2.Range("B2").Value = "TB-16511687079928901.csv"
$ws2.Range("B3").Value = "ZB-match_5-16511687043364024.csv"
$ws2.Range("B4").Value = "ZB-match_0-16511687044504454.csv"
$ws2.Range("B5").Value = "TB-16511687052346437.csv"
$ws2.Range("B6").Value = "OB-16511687048014007.csv"
$ws2.Range("B7").Value = "OB-16511687046713982.csv"
$ws2.Range("B8").Value = "ZB-match_0-16511687046404355.csv"
$ws2.Range("B9").Value = "OB-1651168704757401.csv"
$ws2.Range("B10").Value = "TB-16511687083531144.csv"

# --- Sheet 3: RS_TO (name only, no data changes) -------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651168708377116"

# --- Sheet 4: TOL_TO -------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16511687084401152"
$ws4.Range("B2").Value = "MM_stims-16511687083911135.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687083791122.csv"
$ws4.Range("B4").Value = "MM_stims-16511687084231167.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687083921113.csv"
$ws4.Range("B6").Value = "MM_stims-16511687084391499.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687084241102.csv"

# --- Sheet 5: vSAT_TO -------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511687085341108"
$ws5.Range("B2").Value = "vSAT_stims-16511687085181093.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687084721174.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687084451208.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511687084861488.csv"
